$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(5, 6).Value = 515
$ws.Cells.Item(6, 6).Value = 460
$ws.Cells.Item(9, 6).Value = 16
$ws.Cells.Item(10, 6).Value = 13003
$ws.Cells.Item(11, 6).Value = 13003
$ws.Cells.Item(16, 6).Value = 209
$ws.Cells.Item(17, 6).Value = 178
$ws.Cells.Item(18, 6).Value = 222
$ws.Cells.Item(19, 6).Value = 2800
$ws.Cells.Item(21, 6).Value = 112
$ws.Cells.Item(22, 6).Value = 2144
$ws.Cells.Item(23, 6).Value = 193
$ws.Cells.Item(26, 6).Value = 75
$ws.Cells.Item(27, 6).Value = 2493
$ws.Cells.Item(28, 6).Value = 86
$ws.Cells.Item(29, 6).Value = 1213
$ws.Cells.Item(30, 6).Value = 4457
$ws.Cells.Item(32, 6).Value = 4103
$ws.Cells.Item(33, 6).Value = 1145
$ws.Cells.Item(34, 6).Value = 2733
$ws.Cells.Item(35, 6).Value = 3142
$ws.Cells.Item(36, 6).Value = 117
$ws.Cells.Item(37, 6).Value = 1457
$ws.Cells.Item(39, 6).Value = 812
$ws.Cells.Item(40, 6).Value = 78
$ws.Cells.Item(41, 6).Value = 215
$ws.Cells.Item(42, 6).Value = 756
$ws.Cells.Item(43, 6).Value = 1212
$ws.Cells.Item(44, 6).Value = 104
$ws.Cells.Item(45, 6).Value = 202
$ws.Cells.Item(46, 6).Value = 496
$ws.Cells.Item(47, 6).Value = 143
$ws.Cells.Item(48, 6).Value = 251
$ws.Cells.Item(49, 6).Value = 301

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(3, 6).Value = 81
$ws.Cells.Item(6, 6).Value = 47
$ws.Cells.Item(11, 6).Value = 192
$ws.Cells.Item(13, 6).Value = 21
$ws.Cells.Item(16, 6).Value = 48

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(4, 6).Value = 515
$ws.Cells.Item(5, 6).Value = 460
$ws.Cells.Item(8, 6).Value = 13003
$ws.Cells.Item(11, 6).Value = 81
$ws.Cells.Item(12, 6).Value = 47
$ws.Cells.Item(13, 6).Value = 209
$ws.Cells.Item(14, 6).Value = 178
$ws.Cells.Item(16, 6).Value = 222
$ws.Cells.Item(17, 6).Value = 2800
$ws.Cells.Item(18, 6).Value = 2144
$ws.Cells.Item(19, 6).Value = 193
$ws.Cells.Item(22, 6).Value = 75
$ws.Cells.Item(24, 6).Value = 2493
$ws.Cells.Item(25, 6).Value = 1214
$ws.Cells.Item(26, 6).Value = 192
$ws.Cells.Item(27, 6).Value = 21
$ws.Cells.Item(28, 6).Value = 4457
$ws.Cells.Item(30, 6).Value = 4103
$ws.Cells.Item(31, 6).Value = 1146
$ws.Cells.Item(32, 6).Value = 2733
$ws.Cells.Item(33, 6).Value = 3142
$ws.Cells.Item(34, 6).Value = 117
$ws.Cells.Item(36, 6).Value = 1457
$ws.Cells.Item(37, 6).Value = 48
$ws.Cells.Item(39, 6).Value = 812
$ws.Cells.Item(40, 6).Value = 78
$ws.Cells.Item(41, 6).Value = 215
$ws.Cells.Item(42, 6).Value = 756
$ws.Cells.Item(44, 6).Value = 1212
$ws.Cells.Item(45, 6).Value = 104
$ws.Cells.Item(46, 6).Value = 202
$ws.Cells.Item(47, 6).Value = 496
$ws.Cells.Item(48, 6).Value = 143
$ws.Cells.Item(49, 6).Value = 251
$ws.Cells.Item(50, 6).Value = 301
